# Update "想去人数" (F column) counts on the "展览" and "全部类型" sheets
# to reflect the latest generated data snapshot (gh-pages output @ 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 3494
$ws1.Range("F6").Value = 425
$ws1.Range("F8").Value = 62
$ws1.Range("F9").Value = 49
$ws1.Range("F10").Value = 1289
$ws1.Range("F12").Value = 1713

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 3494
$ws4.Range("F6").Value = 425
$ws4.Range("F9").Value = 62
$ws4.Range("F10").Value = 49
$ws4.Range("F13").Value = 1289
$ws4.Range("F15").Value = 1713
